# Refresh the cryptos list (price / 1h-volume columns) with the latest
# scraped figures, and swap the InjectiveProtocol/Mantle rows (48/49)
# whose relative ranking flipped in this run.
#
# Note: several "Price" values are plain decimals (e.g. "563.45") that
# Excel would otherwise auto-convert to numbers on assignment, losing the
# original text formatting/trailing zeros. A leading apostrophe forces
# those specific cells to stay text, matching the source sheet where every
# Price/Volume cell is stored as a string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.069.60'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '2.422.72'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'563.45"
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = "'143.82"
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = '2.422.16'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -0.59%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  -3.54%  '
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').Value = "'26.19"
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('E15').Value = '  -2.42%  '
$ws.Range('D17').Value = '61.941.85'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '2.442.42'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('D21').Value = "'6.83"
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('E23').Value = '  -0.10%  '
$ws.Range('D24').Value = "'67.29"
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('E25').Value = '  +0.47%  '
$ws.Range('D26').Value = "'8.82"
$ws.Range('E26').Value = '  -2.01%  '
$ws.Range('D27').Value = "'557.39"
$ws.Range('E27').Value = '  -5.56%  '
$ws.Range('D28').Value = '2.542.42'
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').Value = "'1.00"
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('E32').Value = '  -5.64%  '
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('E34').Value = '  -0.82%  '
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').Value = "'4.75"
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range('E39').Value = '  -5.11%  '
$ws.Range('D40').Value = "'152.32"
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').Value = "'18.67"
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('E42').Value = '  -1.57%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  -3.83%  '
$ws.Range('D45').Value = "'147.57"
$ws.Range('E45').Value = '  -1.84%  '
$ws.Range('E46').Value = '  -0.53%  '
$ws.Range('D47').Value = "'0.0531"
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = "'19.93"
$ws.Range('E48').Value = '  -2.65%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.596"
$ws.Range('E49').Value = '  +0.47%  '
$ws.Range('D50').Value = "'0.0920"
$ws.Range('E50').Value = '  -0.74%  '
$ws.Range('E51').Value = '  -0.55%  '
